# Update cryptos list (Price / Volume(1h) columns) with freshly scraped values.
# Note: D-column cells whose new text would otherwise be auto-parsed by Excel
# as a plain number (e.g. "1.00", "590.51") are written with a leading
# apostrophe so they stay plain text, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.114.70'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '2.929.30'
$ws.Range("E3").Value = '  +1.34%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''590.51'
$ws.Range("D6").Value = '''145.52'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.506'
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("D9").Value = '2.927.95'
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").Value = '''6.85'
$ws.Range("E10").Value = '  +3.12%  '
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("E13").Value = '  +1.82%  '
$ws.Range("D14").Value = '''33.73'
$ws.Range("E14").Value = '  -1.10%  '
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("D16").Value = '3.412.90'
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").Value = '61.070.14'
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("D19").Value = '2.923.57'
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("D20").Value = '''431.35'
$ws.Range("E20").Value = '  +1.74%  '
$ws.Range("D21").Value = '''13.48'
$ws.Range("E21").Value = '  -0.65%  '
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("D23").Value = '''7.10'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '''81.19'
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").Value = '''11.11'
$ws.Range("E25").Value = '  +0.66%  '
$ws.Range("E26").Value = '  +2.91%  '
$ws.Range("E27").Value = '  +2.74%  '
$ws.Range("E29").Value = '  +6.89%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").Value = '''2.62'
$ws.Range("E31").Value = '  +0.59%  '
$ws.Range("D32").Value = '''7.11'
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("D33").Value = '''26.53'
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("E34").Value = '  +2.66%  '
$ws.Range("D35").Value = '0.0₃0863'
$ws.Range("E35").Value = '  +3.58%  '
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("D37").Value = '''5.64'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("E38").Value = '  +5.28%  '
$ws.Range("D39").Value = '''49.99'
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("E40").Value = '  +1.76%  '
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").Value = '''8.61'
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = '''39.31'
$ws.Range("E44").Value = '  -4.49%  '
$ws.Range("D45").Value = '''376.41'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("D47").Value = '2.711.31'
$ws.Range("E47").Value = '  +2.54%  '
$ws.Range("D48").Value = '''131.63'
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("D50").Value = '''24.25'
$ws.Range("E50").Value = '  -3.68%  '
$ws.Range("E51").Value = '  +0.58%  '
